$wb = $excel.ActiveWorkbook

# "Sheet1" (the old site-layout / map scratch sheet) is repurposed into a
# new "users" sheet that just shows a preview name. Clear its previous
# content and replace it with the single message/preview cell.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "users"
$ws.Cells.Clear()

$ws.Range("A1").Value = "Eemeli"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

# Relocate "users" so it sits right before "translations" (sites, openings,
# users, translations) and make it the active tab/selection.
$translations = $wb.Worksheets.Item("translations")
$ws.Move($translations)

$users = $wb.Worksheets.Item("users")
$users.Activate()
[void]$users.Range("A1").Select()
